# Add sample EMN001 to dataset S4b
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("S4b - ref. alignment ag. EMN001")

# Insert a new row 5 (shifts existing rows 5-70 down to 6-71),
# preserving formatting from the row below.
$ws.Rows.Item(5).Insert()

$ws.Cells.Item(5, 1).Value = "EMN001"
$ws.Cells.Item(5, 2).Value = 6993330
$ws.Cells.Item(5, 3).Value = 189466499
$ws.Cells.Item(5, 4).Value = 3.69
